$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix inconsistent casing of the class-name labels in column A
$ws.Range("A2").Value = "mdaTextHomePage"
$ws.Range("A8").Value = "pageTitleNewTab"
$ws.Range("A4").Value = "mdaTitle"

# Move the active selection to A4
$ws.Range("A4").Select()
